$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# 1. Update the "Speichern" button xpath selector text in F2
$ws.Range("F2").Value = "xpath=//div[@id='createAbsenceModal']//button[@type='submit' and text()='Speichern']"

# 2. Widen column F so the longer text keeps its "best fit" look
#    (69 chars is the input that lands closest to the recorded 69.77734375
#    once the host quantises it to whole pixels)
$ws.Columns.Item(6).ColumnWidth = 69

# 3. Resize the picture on the sheet: keep the same top-left anchor but
#    shrink it so the right edge lands at column 7 + 727602 EMU instead of
#    column 11 + 331362 EMU (height/bottom unchanged).
$shape = $ws.Shapes.Item(1)
$emuPerPoint = 12700
$newRightEmu = 0
for ($c = 0; $c -lt 7; $c++) {
    $newRightEmu += $ws.Columns.Item($c + 1).Width * $emuPerPoint
}
$newRightEmu += 727602
$newRightPt = $newRightEmu / $emuPerPoint
$shape.Width = $newRightPt - $shape.Left

# 4. Update the selected cell shown when the workbook is opened
$ws.Range("F2").Select()
